$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/96cc2da3e2b22021de61ffcfe3c1a86cca75053f/e2e/ad1403bb-45aa-47ac-aefe-49a40bd189e3.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/c3022a67ef3423dea50f0e43aa7291584603c368/e2e/ad1403bb-45aa-47ac-aefe-49a40bd189e3.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/c3022a67ef3423dea50f0e43aa7291584603c368/e2e/ad1403bb-45aa-47ac-aefe-49a40bd189e3.md"
$displayName = "ad1403bb-45aa-47ac-aefe-49a40bd189e3.md"

# Hyperlink font color used by the existing "HyperLink" style (RGB 6495ED),
# passed BGR-packed the way the COM Color property expects it.
$hyperlinkColor = 15570276

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsZh.Range("I8").Value = $displayName
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestUrl, "", "", $displayName)
$wsZh.Range("I8").Font.Underline = $true
$wsZh.Range("I8").Font.Color = $hyperlinkColor

$wsZh.Range("J8").Value = "ad1403bb-45aa-47ac-aefe-49a40bd189e3.e5e15d44b12b411ac58032bf20780e67b84e837a.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-14 03:00:29"
$wsZh.Range("P8").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDe.Range("I8").Value = $displayName
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestUrl, "", "", $displayName)
$wsDe.Range("I8").Font.Underline = $true
$wsDe.Range("I8").Font.Color = $hyperlinkColor

$wsDe.Range("J8").Value = "ad1403bb-45aa-47ac-aefe-49a40bd189e3.e5e15d44b12b411ac58032bf20780e67b84e837a.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-14 03:00:40"
$wsDe.Range("P8").Value = $errorDetail
